$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previously-scraped data only grabbed team statistics, not the season
# record (Wins/Losses/Ties). Add three new columns (AD:AF) with that record.

# Copy the formatting of the existing last header cell (AC1) into the three
# new header cells so they match the rest of the header row (bold, border,
# centered).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record for every player row (2 through 44): 94 wins, 67 losses, 0 ties
for ($i = 2; $i -le 44; $i++) {
    $ws.Cells.Item($i, 30).Value = 94
    $ws.Cells.Item($i, 31).Value = 67
    $ws.Cells.Item($i, 32).Value = 0
}
